$wb = $excel.ActiveWorkbook

# --- Admin/User flow: add a new login row to the "Details" sheet ---
$details = $wb.Worksheets.Item("Details")
$details.Range("A4").Value = "sru"
$details.Range("B4").Value = "1234"

# --- Reduce WALLE's first showtime capacity from 20 to 18 ---
$walle = $wb.Worksheets.Item("WALLE")
$walle.Range("M2").Value = 18
